# Apply updated cryptocurrency price/volume data (and a couple of
# re-ranked rows where the coin name/link/price/volume moved to a new row)
# to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.743.94'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').Value = '1.886.29'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('D4').Formula = '="0.9996"'
$ws.Range('D4').Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Formula = '="0.7933"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  -5.61%  '
$ws.Range('D6').Formula = '="241.36"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('D7').Formula = '="0.9998"'
$ws.Range('D7').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Formula = '="0.3172"'
$ws.Range('D8').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  -1.98%  '
$ws.Range('E9').Value = '  -4.86%  '
$ws.Range('D10').Formula = '="0.07003"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('D11').Formula = '="0.08034"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').Formula = '="0.7616"'
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  +1.13%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.901.67'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Formula = '="5.296"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  +1.29%  '
$ws.Range('D15').Formula = '="92.23"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').Value = '29.696.66'
$ws.Range('E16').Value = '  -1.12%  '
$ws.Range('D17').Formula = '="13.85"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  -2.52%  '
$ws.Range('D18').Formula = '="5.930"'
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  -0.59%  '
$ws.Range('D19').Formula = '="243.38"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  -0.69%  '
$ws.Range('D20').Formula = '="0.000007681"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  -1.25%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').Formula = '="8.220"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  +17.23%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Formula = '="0.9995"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = '2.126.73'
$ws.Range('E23').Value = '  -1.24%  '
$ws.Range('D24').Formula = '="0.9999"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Formula = '="0.1683"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +3.97%  '
$ws.Range('D26').Formula = '="9.290"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').Formula = '="164.22"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  -3.20%  '
$ws.Range('D28').Formula = '="18.61"'
$ws.Range('D28').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  -1.92%  '
$ws.Range('D29').Formula = '="2.050"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  -1.75%  '
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').Formula = '="1.532"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  +0.97%  '
$ws.Range('D32').Formula = '="4.376"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('D33').Formula = '="0.05675"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  +0.38%  '
$ws.Range('D34').Formula = '="4.048"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('D35').Formula = '="1.261"'
$ws.Range('D35').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  -2.02%  '
$ws.Range('D36').Formula = '="0.7340"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  -0.40%  '
$ws.Range('D37').Formula = '="0.9954"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('D38').Formula = '="2.615"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  -3.83%  '
$ws.Range('D39').Formula = '="0.01907"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('D40').Formula = '="2.767"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('D41').Formula = '="0.4403"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  -1.02%  '
$ws.Range('D42').Formula = '="72.43"'
$ws.Range('D42').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('D43').Formula = '="5.815"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  -3.44%  '
$ws.Range('D44').Formula = '="0.9994"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').Formula = '="0.8351"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('D46').Formula = '="102.55"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +1.13%  '
$ws.Range('D47').Value = '1.019.99'
$ws.Range('E47').Value = '  +3.19%  '
$ws.Range('D48').Formula = '="1.865"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  -2.15%  '
$ws.Range('D49').Formula = '="9.864"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  +0.90%  '
$ws.Range('D50').Formula = '="7.417"'
$ws.Range('D50').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  -2.74%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Formula = '="2.901"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +5.38%  '

$excel.CutCopyMode = 0

